# "Mejora para obtener el df_fatiga_media"
# Adds a new "Sheet2" after "Sheet1" holding the computed fatigue-percentage
# series: row index (A), date (B) and "Porcentaje de fatiga" (C).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so it lands as the 2nd tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Formats ------------------------------------------------------------
# Reuse Sheet1's existing cell formats (bold/centered index column, date
# column) instead of re-describing them, so no new style entries get minted.
$ws1.Range("A2").Copy()
$ws2.Range("A2:A11").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("C2").Copy()
$ws2.Range("B2:B11").PasteSpecial(-4122)   # xlPasteFormats (date format)

$ws1.Range("B1").Copy()
$ws2.Range("B1:C1").PasteSpecial(-4122)    # xlPasteFormats (header format)

$excel.CutCopyMode = $false

# --- Header row -----------------------------------------------------------
$ws2.Range("B1").Value = "Fecha"
$ws2.Range("C1").Value = "Porcentaje de fatiga"

# --- Data -----------------------------------------------------------------
$idx    = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9)
$fecha  = @(44956, 44964, 44966, 44971, 44976, 44978, 44980, 44984, 44986, 44991)
$fatiga = @(
    13.28449328449329,
    12.57682177348552,
    20.94769321187188,
    17.34972677595628,
    6.451612903225811,
    18.64406779661016,
    18.57683982683983,
    14.3905360886493,
    16.12903225806451,
    21.66832834828847
)

for ($i = 0; $i -lt $idx.Count; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $idx[$i]
    $ws2.Range("B$r").Value = $fecha[$i]
    $ws2.Range("C$r").Value = $fatiga[$i]
}

# Restore Sheet1 as the active/selected tab (unchanged from before the edit).
$ws1.Activate()

Write-Output "Sheet2 added with $($idx.Count) fatigue rows"
